$d = $word.ActiveDocument

# 1. Update the date field result: "25. prosince 2020" -> "29. prosince 2020"
[void]$d.Content.Find.Execute("25. prosince 2020", $true, $false, $false, $false, $false,
                         $true, 1, $false, "29. prosince 2020", 2)

# 2. Locate the "Tabulka souborů" heading paragraph (currently holds
#    <w:lastRenderedPageBreak/> + the heading text) and retask it as the
#    new "Použité technologie" heading, then insert the new technology
#    paragraph, a blank spacer paragraph, and a fresh "Tabulka souborů"
#    heading paragraph after it (in that order), right before the table.
$heading = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.TrimEnd([char]13, [char]7) -eq "Tabulka souborů") {
        $heading = $cand
        break
    }
}

$heading.Range.Text = "Použité technologie"

$r = $heading.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$techPara = $heading.Next()
$techPara.Style = "Normal"
$techPara.Range.Text = "Zvoleným programovacím jazykem pro vytvoření tohoto projektu je C# 9. Byla vytvořena desktopová aplikace nad frameworkem Windows Forms a platformou .NET 5. Použitými knihovnami jsou NAudio (načtení vzorků z .wav souborů) a OxyPlot (grafy)."

$r2 = $techPara.Range
$r2.Collapse(0)
$r2.InsertParagraphAfter()
$blankPara = $techPara.Next()
$blankPara.Style = "Normal"

$r3 = $blankPara.Range
$r3.Collapse(0)
$r3.InsertParagraphAfter()
$tablePara = $blankPara.Next()
$tablePara.Style = "Heading1"
$tablePara.Range.Text = "Tabulka souborů"

# 3. At the end of the document, turn the final empty paragraph into the
#    "Rozdělení na rámce" heading and append the frame-splitting paragraph
#    after it.
$lastPara = $d.Paragraphs.Last
$lastPara.Style = "Heading1"
$lastPara.Range.Text = "Rozdělení na rámce"

$r4 = $lastPara.Range
$r4.Collapse(0)
$r4.InsertParagraphAfter()
$framePara = $d.Paragraphs.Last
$framePara.Style = "Normal"
$framePara.Range.Text = "Extrahovaná 1s nahrávky je ustředněna, normalizována a rozdělena na rámce o délce 20ms, překrývající se po 10ms. Celkem tedy na 1s máme 99 rámců."
